$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 56002
$ws.Cells.Item(2, 4).Value = 113796244
$ws.Cells.Item(3, 3).Value = 135912
$ws.Cells.Item(3, 4).Value = 318240905
$ws.Cells.Item(4, 3).Value = 49318
$ws.Cells.Item(4, 4).Value = 142484849
$ws.Cells.Item(5, 3).Value = 15570
$ws.Cells.Item(5, 4).Value = 52268166
$ws.Cells.Item(6, 3).Value = 5711
$ws.Cells.Item(6, 4).Value = 25871245
$ws.Cells.Item(7, 3).Value = 1120
$ws.Cells.Item(7, 4).Value = 6479790
$ws.Cells.Item(8, 3).Value = 66
$ws.Cells.Item(8, 4).Value = 449315
$ws.Cells.Item(12, 3).Value = 58170
$ws.Cells.Item(12, 4).Value = 92355058
$ws.Cells.Item(13, 3).Value = 14209
$ws.Cells.Item(13, 4).Value = 28550055
$ws.Cells.Item(14, 3).Value = 38282
$ws.Cells.Item(14, 4).Value = 88142006
$ws.Cells.Item(15, 3).Value = 12708
$ws.Cells.Item(15, 4).Value = 34991682
$ws.Cells.Item(16, 3).Value = 3658
$ws.Cells.Item(16, 4).Value = 11169437
$ws.Cells.Item(17, 3).Value = 1183
$ws.Cells.Item(17, 4).Value = 5065626
$ws.Cells.Item(20, 3).Value = 14285
$ws.Cells.Item(20, 4).Value = 22126736
$ws.Cells.Item(21, 3).Value = 19994
$ws.Cells.Item(21, 4).Value = 42189883
$ws.Cells.Item(22, 3).Value = 47365
$ws.Cells.Item(22, 4).Value = 113498100
$ws.Cells.Item(23, 3).Value = 16415
$ws.Cells.Item(23, 4).Value = 47193289
$ws.Cells.Item(24, 3).Value = 4849
$ws.Cells.Item(24, 4).Value = 15614940
$ws.Cells.Item(25, 3).Value = 1547
$ws.Cells.Item(25, 4).Value = 6408894
$ws.Cells.Item(26, 3).Value = 246
$ws.Cells.Item(26, 4).Value = 1339889
$ws.Cells.Item(28, 3).Value = 15988
$ws.Cells.Item(28, 4).Value = 24628691
$ws.Cells.Item(29, 3).Value = 11290
$ws.Cells.Item(29, 4).Value = 22993510
$ws.Cells.Item(30, 3).Value = 32744
$ws.Cells.Item(30, 4).Value = 73770030
$ws.Cells.Item(31, 3).Value = 11797
$ws.Cells.Item(31, 4).Value = 31718147
$ws.Cells.Item(32, 3).Value = 3233
$ws.Cells.Item(32, 4).Value = 9617169
$ws.Cells.Item(33, 3).Value = 1008
$ws.Cells.Item(33, 4).Value = 4232029
$ws.Cells.Item(36, 3).Value = 11561
$ws.Cells.Item(36, 4).Value = 17937026
$ws.Cells.Item(37, 3).Value = 4995
$ws.Cells.Item(37, 4).Value = 10948734
$ws.Cells.Item(38, 3).Value = 11747
$ws.Cells.Item(38, 4).Value = 27196326
$ws.Cells.Item(39, 3).Value = 4859
$ws.Cells.Item(39, 4).Value = 13780057
$ws.Cells.Item(41, 3).Value = 426
$ws.Cells.Item(41, 4).Value = 2024852
$ws.Cells.Item(44, 3).Value = 3552
$ws.Cells.Item(44, 4).Value = 5451898
$ws.Cells.Item(45, 3).Value = 25564
$ws.Cells.Item(45, 4).Value = 52661676
$ws.Cells.Item(46, 3).Value = 75944
$ws.Cells.Item(46, 4).Value = 178543528
$ws.Cells.Item(47, 3).Value = 29197
$ws.Cells.Item(47, 4).Value = 81003307
$ws.Cells.Item(48, 3).Value = 9502
$ws.Cells.Item(48, 4).Value = 28825959
$ws.Cells.Item(49, 3).Value = 3261
$ws.Cells.Item(49, 4).Value = 13081208
$ws.Cells.Item(51, 3).Value = 29
$ws.Cells.Item(51, 4).Value = 141461
$ws.Cells.Item(53, 3).Value = 26078
$ws.Cells.Item(53, 4).Value = 47516039
$ws.Cells.Item(54, 3).Value = 2705
$ws.Cells.Item(54, 4).Value = 4382992
$ws.Cells.Item(55, 3).Value = 9019
$ws.Cells.Item(55, 4).Value = 14916274
$ws.Cells.Item(56, 3).Value = 3027
$ws.Cells.Item(56, 4).Value = 5328937
$ws.Cells.Item(57, 3).Value = 990
$ws.Cells.Item(57, 4).Value = 1934683
$ws.Cells.Item(58, 3).Value = 299
$ws.Cells.Item(58, 4).Value = 657437
$ws.Cells.Item(59, 3).Value = 44
$ws.Cells.Item(59, 4).Value = 161660
$ws.Cells.Item(61, 3).Value = 9205
$ws.Cells.Item(61, 4).Value = 13642626
$ws.Cells.Item(62, 3).Value = 1821
$ws.Cells.Item(62, 4).Value = 3969477
$ws.Cells.Item(63, 3).Value = 4310
$ws.Cells.Item(63, 4).Value = 9336706
$ws.Cells.Item(64, 3).Value = 1722
$ws.Cells.Item(64, 4).Value = 3867459
$ws.Cells.Item(68, 3).Value = 2819
$ws.Cells.Item(68, 4).Value = 5563831
$ws.Cells.Item(69, 3).Value = 22811
$ws.Cells.Item(69, 4).Value = 45012298
$ws.Cells.Item(70, 3).Value = 66172
$ws.Cells.Item(70, 4).Value = 150475457
$ws.Cells.Item(71, 3).Value = 24178
$ws.Cells.Item(71, 4).Value = 66794992
$ws.Cells.Item(72, 3).Value = 7558
$ws.Cells.Item(72, 4).Value = 22734275
$ws.Cells.Item(73, 3).Value = 2440
$ws.Cells.Item(73, 4).Value = 9830524
$ws.Cells.Item(74, 3).Value = 482
$ws.Cells.Item(74, 4).Value = 2662413
$ws.Cells.Item(78, 3).Value = 21145
$ws.Cells.Item(78, 4).Value = 32432994
$ws.Cells.Item(79, 3).Value = 83056
$ws.Cells.Item(79, 4).Value = 170053249
$ws.Cells.Item(80, 3).Value = 225845
$ws.Cells.Item(80, 4).Value = 507450135
$ws.Cells.Item(81, 3).Value = 101902
$ws.Cells.Item(81, 4).Value = 284155122
$ws.Cells.Item(82, 3).Value = 37119
$ws.Cells.Item(82, 4).Value = 124137090
$ws.Cells.Item(83, 3).Value = 13612
$ws.Cells.Item(83, 4).Value = 61137742
$ws.Cells.Item(84, 3).Value = 2624
$ws.Cells.Item(84, 4).Value = 16749176
$ws.Cells.Item(89, 3).Value = 8
$ws.Cells.Item(89, 4).Value = 28424
$ws.Cells.Item(90, 3).Value = 78941
$ws.Cells.Item(90, 4).Value = 125194999
$ws.Cells.Item(91, 3).Value = 5617
$ws.Cells.Item(91, 4).Value = 8736894
$ws.Cells.Item(92, 3).Value = 13542
$ws.Cells.Item(92, 4).Value = 21409597
$ws.Cells.Item(93, 3).Value = 4351
$ws.Cells.Item(93, 4).Value = 7033128
$ws.Cells.Item(98, 3).Value = 6350
$ws.Cells.Item(98, 4).Value = 8739136
$ws.Cells.Item(100, 3).Value = 7304
$ws.Cells.Item(100, 4).Value = 14507948
$ws.Cells.Item(101, 3).Value = 2623
$ws.Cells.Item(101, 4).Value = 5918186
$ws.Cells.Item(102, 3).Value = 967
$ws.Cells.Item(102, 4).Value = 2293966
$ws.Cells.Item(104, 3).Value = 64
$ws.Cells.Item(104, 4).Value = 283559
$ws.Cells.Item(106, 3).Value = 4889
$ws.Cells.Item(106, 4).Value = 7112333
$ws.Cells.Item(107, 3).Value = 1033
$ws.Cells.Item(107, 4).Value = 2332697
$ws.Cells.Item(109, 3).Value = 253
$ws.Cells.Item(109, 4).Value = 572692
$ws.Cells.Item(113, 3).Value = 16353
$ws.Cells.Item(113, 4).Value = 34280846
$ws.Cells.Item(114, 3).Value = 43061
$ws.Cells.Item(114, 4).Value = 100653899
$ws.Cells.Item(115, 3).Value = 15247
$ws.Cells.Item(115, 4).Value = 42461385
$ws.Cells.Item(116, 3).Value = 4700
$ws.Cells.Item(116, 4).Value = 14782773
$ws.Cells.Item(117, 3).Value = 1456
$ws.Cells.Item(117, 4).Value = 6152150
$ws.Cells.Item(122, 3).Value = 13483
$ws.Cells.Item(122, 4).Value = 20647858
$ws.Cells.Item(123, 3).Value = 44544
$ws.Cells.Item(123, 4).Value = 90057390
$ws.Cells.Item(124, 3).Value = 95385
$ws.Cells.Item(124, 4).Value = 214146873
$ws.Cells.Item(125, 3).Value = 32232
$ws.Cells.Item(125, 4).Value = 86353174
$ws.Cells.Item(126, 3).Value = 10077
$ws.Cells.Item(126, 4).Value = 30443822
$ws.Cells.Item(127, 3).Value = 3184
$ws.Cells.Item(127, 4).Value = 13011005
$ws.Cells.Item(129, 3).Value = 31
$ws.Cells.Item(129, 4).Value = 155332
$ws.Cells.Item(132, 3).Value = 34937
$ws.Cells.Item(132, 4).Value = 53397547
$ws.Cells.Item(133, 3).Value = 53619
$ws.Cells.Item(133, 4).Value = 109767230
$ws.Cells.Item(134, 3).Value = 112284
$ws.Cells.Item(134, 4).Value = 250235098
$ws.Cells.Item(135, 3).Value = 36366
$ws.Cells.Item(135, 4).Value = 99714679
$ws.Cells.Item(136, 3).Value = 10762
$ws.Cells.Item(136, 4).Value = 33134996
$ws.Cells.Item(137, 3).Value = 3354
$ws.Cells.Item(137, 4).Value = 13725261
$ws.Cells.Item(142, 3).Value = 43634
$ws.Cells.Item(142, 4).Value = 65330926
$ws.Cells.Item(143, 3).Value = 19568
$ws.Cells.Item(143, 4).Value = 40030026
$ws.Cells.Item(144, 3).Value = 47604
$ws.Cells.Item(144, 4).Value = 111607451
$ws.Cells.Item(145, 3).Value = 17882
$ws.Cells.Item(145, 4).Value = 49815975
$ws.Cells.Item(146, 3).Value = 5137
$ws.Cells.Item(146, 4).Value = 15774573
$ws.Cells.Item(147, 3).Value = 1525
$ws.Cells.Item(147, 4).Value = 6418903
$ws.Cells.Item(148, 3).Value = 339
$ws.Cells.Item(148, 4).Value = 1969602
$ws.Cells.Item(152, 3).Value = 14751
$ws.Cells.Item(152, 4).Value = 22827878
$ws.Cells.Item(153, 3).Value = 53130
$ws.Cells.Item(153, 4).Value = 109773820
$ws.Cells.Item(154, 3).Value = 123474
$ws.Cells.Item(154, 4).Value = 284207581
$ws.Cells.Item(155, 3).Value = 39279
$ws.Cells.Item(155, 4).Value = 112391605
$ws.Cells.Item(156, 3).Value = 11723
$ws.Cells.Item(156, 4).Value = 39250228
$ws.Cells.Item(157, 3).Value = 4196
$ws.Cells.Item(157, 4).Value = 18741434
$ws.Cells.Item(158, 3).Value = 848
$ws.Cells.Item(158, 4).Value = 5153371
$ws.Cells.Item(160, 3).Value = 40908
$ws.Cells.Item(160, 4).Value = 63528500
